$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Public_Schools")

# Delete the three rows whose CITY is "Hyde Park" (rows 2, 11, 14),
# working from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(2).Delete()

$ws.Range("A12:XFD12").Select()
